# Update the "想去人数" (F column) counts on both the "展览" and "全部类型"
# worksheets to reflect newly generated output data.

$wb = $excel.ActiveWorkbook

# Row -> new F-column value (identical update applied to both sheets).
$updates = @{
    2  = 1196
    3  = 985
    4  = 298
    5  = 65
    8  = 2435
    9  = 7969
    10 = 946
    11 = 483
    12 = 427
    13 = 192
    14 = 451
    16 = 174
    17 = 8241
    18 = 331
    19 = 1425
    24 = 353
    25 = 196
    28 = 121
    29 = 38
    30 = 437
    31 = 1178
    32 = 46
    33 = 62
    34 = 108
    37 = 47
    38 = 84
    39 = 76
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
